# Automatische test-sync: 2025-06-22 18:46:50
# Append a new incoming-mail log entry to the "Logs" sheet and refresh the
# "Dashboard" category summary + conditional-formatting ranges to match.

$wb  = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new row (row 16) to the Logs sheet -----------------
$logs.Range("A16").Value = "Vragen over handleiding"
$logs.Range("B16").Value = "mailmind.test@zohomail.eu"
$logs.Range("C16").Value = "Waar kan ik de handleiding van product X vinden?"
$logs.Range("D16").Value = "Productinformatie"
$logs.Range("E16").Value = "Beste klant,
Bedankt voor uw vraag. De handleiding van product X is te vinden op onze website onder de sectie 'Support' of 'Downloads'. Mocht u hier toch problemen mee ervaren, laat het ons dan weten zodat we u verder kunnen helpen.
Met vriendelijke groet,
[Bedrijfsnaam]"
$logs.Range("F16").Value = "2025-06-22 18:46:33"
$logs.Range("G16").Value = "Ja"

# Keep the row height at the default (no explicit row height), matching
# the rest of the sheet, instead of the auto row-height Excel would apply
# after writing wrapped multi-line text.
$logs.Rows.Item(16).AutoFit()

# --- 2. Extend the conditional formatting ranges to include row 16 ----
$logs.Range("D2:D15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D16"))
$logs.Range("G2:G15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G16"))

# --- 3. Refresh the Dashboard category counts --------------------------
# "Productinformatie" now has 3 occurrences (was 2) and overtakes
# "Sollicitatie / Vacature" (still 2) in the ranking, so the two rows swap.
$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 3
$dash.Range("A4").Value = "Sollicitatie / Vacature"
$dash.Range("B4").Value = 2
